$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.894.72"
$ws.Range("D3").Value = "2.092.98"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.49"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0782"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").Value = "2.389.15"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.16"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.767"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.27"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "2.088.13"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "37.816.21"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  -3.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.83"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.60"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.67"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.141"
$ws.Range("E27").Value = "  +11.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.96"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.50"
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.64"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.51"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.43"
$ws.Range("E36").Value = "  +6.13%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.83"
$ws.Range("E37").Value = "  +3.52%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.45"
$ws.Range("E39").Value = "  -4.09%  "
$ws.Range("E40").Value = "  +7.43%  "
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.48"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").Value = "1.451.83"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E46").Value = "  +3.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.71"
$ws.Range("E47").Value = "  +4.19%  "
$ws.Range("E48").Value = "  -8.85%  "
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.03"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").Value = "2.285.07"
$ws.Range("E51").Value = "  +1.06%  "

Write-Host "Applied cryptos update"
